{"js": "const replacements = [\n  [\"2024-05-12 Sunday\", \"2024-05-13 Monday\"],\n  [\"575\u00f75=\", \"594\u00f77=\"],\n  [\"442\u00f73=\", \"750\u00f72=\"],\n  [\"869\u00f76=\", \"799\u00f76=\"],\n  [\"491\u00f73=\", \"554\u00f74=\"],\n  [\"908\u00f72=\", \"725\u00f77=\"],\n  [\"679\u00f74=\", \"988\u00f77=\"],\n  [\"148\u00f79=\", \"234\u00f72=\"],\n  [\"180\u00f78=\", \"440\u00f73=\"],\n  [\"627\u00f77=\", \"823\u00f76=\"],\n  [\"227\u00f74=\", \"728\u00f74=\"],\n  [\"458\u00f72=\", \"571\u00f72=\"],\n  [\"894\u00f79=\", \"420\u00f74=\"],\n  [\"915\u00f75=\", \"746\u00f77=\"],\n  [\"887\u00f79=\", \"671\u00f79=\"],\n  [\"675\u00f75=\", \"720\u00f79=\"],\n  [\"782\u00f73=\", \"788\u00f75=\"],\n  [\"270\u00f75=\", \"867\u00f79=\"],\n  [\"759\u00f74=\", \"963\u00f73=\"],\n  [\"742\u00f75=\", \"753\u00f78=\"],\n  [\"714\u00f78=\", \"793\u00f79=\"],\n  [\"373\u00f72=\", \"848\u00f78=\"],\n  [\"885\u00f73=\", \"542\u00f75=\"],\n  [\"124\u00f76=\", \"308\u00f79=\"],\n  [\"605\u00f73=\", \"668\u00f75=\"],\n  [\"168\u00f77=\", \"175\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2024-05-12 Sunday\", \"2024-05-13 Monday\"),\n  @(\"575\u00f75=\", \"594\u00f77=\"),\n  @(\"442\u00f73=\", \"750\u00f72=\"),\n  @(\"869\u00f76=\", \"799\u00f76=\"),\n  @(\"491\u00f73=\", \"554\u00f74=\"),\n  @(\"908\u00f72=\", \"725\u00f77=\"),\n  @(\"679\u00f74=\", \"988\u00f77=\"),\n  @(\"148\u00f79=\", \"234\u00f72=\"),\n  @(\"180\u00f78=\", \"440\u00f73=\"),\n  @(\"627\u00f77=\", \"823\u00f76=\"),\n  @(\"227\u00f74=\", \"728\u00f74=\"),\n  @(\"458\u00f72=\", \"571\u00f72=\"),\n  @(\"894\u00f79=\", \"420\u00f74=\"),\n  @(\"915\u00f75=\", \"746\u00f77=\"),\n  @(\"887\u00f79=\", \"671\u00f79=\"),\n  @(\"675\u00f75=\", \"720\u00f79=\"),\n  @(\"782\u00f73=\", \"788\u00f75=\"),\n  @(\"270\u00f75=\", \"867\u00f79=\"),\n  @(\"759\u00f74=\", \"963\u00f73=\"),\n  @(\"742\u00f75=\", \"753\u00f78=\"),\n  @(\"714\u00f78=\", \"793\u00f79=\"),\n  @(\"373\u00f72=\", \"848\u00f78=\"),\n  @(\"885\u00f73=\", \"542\u00f75=\"),\n  @(\"124\u00f76=\", \"308\u00f79=\"),\n  @(\"605\u00f73=\", \"668\u00f75=\"),\n  @(\"168\u00f77=\", \"175\u00f76=\")\n)\n\nforeach ($pair in $replacements) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
